$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the trailing "999999 / Clan Level" block of rows (old rows 144-180).
$ws.Range("A144:E180").EntireRow.Delete()

# 2. Remove the old row 133 (its values are being fully replaced below).
$ws.Range("A133:E133").EntireRow.Delete()

# 3. Insert 6 fresh rows at 133 to host the new leaderboard entries
#    (old rows 134.. shift down to make room).
$ws.Range("A133:E138").EntireRow.Insert()

# 4. Populate the new rows with the updated data.
#    A leading "'" forces Excel to keep the numeric-looking IDs/ranks as
#    text (matching the rest of the sheet); .Style restores the default
#    "Normal" look so no stray quote-prefix formatting is left behind.

$ws.Cells.Item(133,1).Value = "'47384"
$ws.Cells.Item(133,1).Style = "Normal"
$ws.Cells.Item(133,2).Value = "'41849539"
$ws.Cells.Item(133,2).Style = "Normal"
$ws.Cells.Item(133,3).Value = "三号馆馆主"
$ws.Cells.Item(133,4).Value = "三馆"
$ws.Cells.Item(133,5).Value = "'3168"
$ws.Cells.Item(133,5).Style = "Normal"

$ws.Cells.Item(134,1).Value = "'0"
$ws.Cells.Item(134,1).Style = "Normal"
$ws.Cells.Item(134,2).Value = "'46248210"
$ws.Cells.Item(134,2).Style = "Normal"
$ws.Cells.Item(134,3).Value = "Ajay"
$ws.Cells.Item(134,4).Value = "三馆"
$ws.Cells.Item(134,5).Value = "'1218"
$ws.Cells.Item(134,5).Style = "Normal"

$ws.Cells.Item(135,1).Value = "'79079"
$ws.Cells.Item(135,1).Style = "Normal"
$ws.Cells.Item(135,2).Value = "'47533851"
$ws.Cells.Item(135,2).Style = "Normal"
$ws.Cells.Item(135,3).Value = "Bibek"
$ws.Cells.Item(135,4).Value = "三馆"
$ws.Cells.Item(135,5).Value = "'2271"
$ws.Cells.Item(135,5).Style = "Normal"

$ws.Cells.Item(136,1).Value = "'47413"
$ws.Cells.Item(136,1).Style = "Normal"
$ws.Cells.Item(136,2).Value = "'47622456"
$ws.Cells.Item(136,2).Style = "Normal"
$ws.Cells.Item(136,3).Value = "伊恩"
$ws.Cells.Item(136,4).Value = "三馆"
$ws.Cells.Item(136,5).Value = "'3166"
$ws.Cells.Item(136,5).Style = "Normal"

$ws.Cells.Item(137,1).Value = "'66453"
$ws.Cells.Item(137,1).Style = "Normal"
$ws.Cells.Item(137,2).Value = "'49553719"
$ws.Cells.Item(137,2).Style = "Normal"
$ws.Cells.Item(137,3).Value = '"Oreo Captain Sir"'
$ws.Cells.Item(137,4).Value = "三馆"
$ws.Cells.Item(137,5).Value = "'2530"
$ws.Cells.Item(137,5).Style = "Normal"

$ws.Cells.Item(138,1).Value = "'0"
$ws.Cells.Item(138,1).Style = "Normal"
$ws.Cells.Item(138,2).Value = "'50742014"
$ws.Cells.Item(138,2).Style = "Normal"
$ws.Cells.Item(138,3).Value = '"Aditya Naik"'
$ws.Cells.Item(138,4).Value = "三馆"
$ws.Cells.Item(138,5).Value = "'997"
$ws.Cells.Item(138,5).Style = "Normal"
